# "Generate Report for Handback"
# The handback transform failed for the f347b540-... file in both the
# zh-cn and de-de locales. Update the Status column (C) for that row to
# reflect the failure, and populate the Error Detail column (P) with the
# specific mismatch message for each locale. Also widen the Error Detail
# column so the message is readable.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# Row 3 everywhere corresponds to the f347b540-... file, whose status text
# ("Ready for handoff") is shared across all three sheets via the shared
# string table. Update every occurrence identically.
$ws1.Range("E3").Value = "Handback transform failed"
$ws1.Range("F3").Value = "Handback transform failed"
$ws2.Range("C3").Value = "Handback transform failed"
$ws3.Range("C3").Value = "Handback transform failed"

# Error Detail (column P) for that same row.
$ws2.Range("P3").Value = "Handback file name: 13gjbken.vda is different with handoff file name: f347b540-b497-4ea3-a6b8-96e45a9cee06.a233eb2ea86bf423053c7516f8bb89a42bf436f0.zh-cn."
$ws3.Range("P3").Value = "Handback file name: 13gjbken.vda is different with handoff file name: f347b540-b497-4ea3-a6b8-96e45a9cee06.a233eb2ea86bf423053c7516f8bb89a42bf436f0.de-de."

# Widen the Error Detail column to fit the new, longer message. The engine
# stores XML column width as ColumnWidth + 0.8333 (5/6, default cell
# padding), so use 39.17 here to land on a stored width of 40.
$ws2.Columns.Item(16).ColumnWidth = 39.17
$ws3.Columns.Item(16).ColumnWidth = 39.17
